function Set-TextValue {
    param($range, [string]$text)
    $origStyle = $range.Style
    $range.NumberFormat = "@"
    $range.Value = $text
    $range.Style = $origStyle
}

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

Set-TextValue $ws.Range('D2') '66.962.55'
Set-TextValue $ws.Range('E2') '  -0.11%  '
Set-TextValue $ws.Range('D3') '3.121.04'
Set-TextValue $ws.Range('E3') '  +0.44%  '
Set-TextValue $ws.Range('E4') '  +0.01%  '
Set-TextValue $ws.Range('D5') '580.49'
Set-TextValue $ws.Range('E5') '  +0.06%  '
Set-TextValue $ws.Range('D6') '173.21'
Set-TextValue $ws.Range('E6') '  +0.45%  '
Set-TextValue $ws.Range('E7') '  -0.05%  '
Set-TextValue $ws.Range('D8') '0.522'
Set-TextValue $ws.Range('E8') '  -0.34%  '
Set-TextValue $ws.Range('D9') '6.43'
Set-TextValue $ws.Range('E9') '  -2.18%  '
Set-TextValue $ws.Range('D10') '0.154'
Set-TextValue $ws.Range('E10') '  -1.23%  '
Set-TextValue $ws.Range('E11') '  -0.66%  '
Set-TextValue $ws.Range('D12') '0.0000248'
Set-TextValue $ws.Range('E12') '  -0.86%  '
Set-TextValue $ws.Range('D13') '37.43'
Set-TextValue $ws.Range('E13') '  +1.23%  '
Set-TextValue $ws.Range('E14') '  -1.51%  '
Set-TextValue $ws.Range('B15') 'WrappedBTC'
Set-TextValue $ws.Range('C15') 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
Set-TextValue $ws.Range('D15') '66.893.63'
Set-TextValue $ws.Range('E15') '  -0.10%  '
Set-TextValue $ws.Range('B16') 'Polkadot'
Set-TextValue $ws.Range('C16') 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
Set-TextValue $ws.Range('D16') '7.14'
Set-TextValue $ws.Range('E16') '  -0.83%  '
Set-TextValue $ws.Range('B17') 'WrappedEther'
Set-TextValue $ws.Range('C17') 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
Set-TextValue $ws.Range('D17') '3.118.90'
Set-TextValue $ws.Range('E17') '  +0.69%  '
Set-TextValue $ws.Range('B18') 'Chainlink'
Set-TextValue $ws.Range('C18') 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
Set-TextValue $ws.Range('D18') '16.42'
Set-TextValue $ws.Range('E18') '  +1.46%  '
Set-TextValue $ws.Range('B19') 'BitcoinCash'
Set-TextValue $ws.Range('C19') 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
Set-TextValue $ws.Range('D19') '485.83'
Set-TextValue $ws.Range('E19') '  +1.34%  '
Set-TextValue $ws.Range('B20') 'Polygon'
Set-TextValue $ws.Range('C20') 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
Set-TextValue $ws.Range('D20') '0.709'
Set-TextValue $ws.Range('E20') '  -0.98%  '
Set-TextValue $ws.Range('B21') 'Uniswap'
Set-TextValue $ws.Range('C21') 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
Set-TextValue $ws.Range('D21') '7.82'
Set-TextValue $ws.Range('E21') '  +4.01%  '
Set-TextValue $ws.Range('B22') 'Litecoin'
Set-TextValue $ws.Range('C22') 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
Set-TextValue $ws.Range('D22') '84.11'
Set-TextValue $ws.Range('E22') '  +0.03%  '
Set-TextValue $ws.Range('B23') 'InternetComputer(DFINITY)'
Set-TextValue $ws.Range('C23') 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
Set-TextValue $ws.Range('D23') '13.21'
Set-TextValue $ws.Range('E23') '  +1.18%  '
Set-TextValue $ws.Range('B24') 'Fetch.AI'
Set-TextValue $ws.Range('C24') 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
Set-TextValue $ws.Range('D24') '2.28'
Set-TextValue $ws.Range('E24') '  -2.25%  '
Set-TextValue $ws.Range('B25') 'RenderToken'
Set-TextValue $ws.Range('C25') 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
Set-TextValue $ws.Range('D25') '10.38'
Set-TextValue $ws.Range('E25') '  +3.08%  '
Set-TextValue $ws.Range('B26') 'Dai'
Set-TextValue $ws.Range('C26') 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
Set-TextValue $ws.Range('D26') '1.00'
Set-TextValue $ws.Range('E26') '  +0.01%  '
Set-TextValue $ws.Range('B27') 'NEARProtocol'
Set-TextValue $ws.Range('C27') 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
Set-TextValue $ws.Range('D27') '7.93'
Set-TextValue $ws.Range('E27') '  -1.21%  '
Set-TextValue $ws.Range('B28') 'ImmutableX'
Set-TextValue $ws.Range('C28') 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
Set-TextValue $ws.Range('D28') '2.36'
Set-TextValue $ws.Range('E28') '  -1.68%  '
Set-TextValue $ws.Range('B29') 'PancakeSwap'
Set-TextValue $ws.Range('C29') 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
Set-TextValue $ws.Range('D29') '2.68'
Set-TextValue $ws.Range('E29') '  +0.19%  '
Set-TextValue $ws.Range('B30') 'EthereumClassic'
Set-TextValue $ws.Range('C30') 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
Set-TextValue $ws.Range('D30') '28.68'
Set-TextValue $ws.Range('E30') '  +1.11%  '
Set-TextValue $ws.Range('B31') 'Hedera'
Set-TextValue $ws.Range('C31') 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
Set-TextValue $ws.Range('D31') '0.115'
Set-TextValue $ws.Range('E31') '  -0.46%  '
Set-TextValue $ws.Range('B32') 'PEPE'
Set-TextValue $ws.Range('C32') 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
Set-TextValue $ws.Range('D32') '0.0₃0952'
Set-TextValue $ws.Range('E32') '  -7.10%  '
Set-TextValue $ws.Range('B33') 'FirstDigitalUSD'
Set-TextValue $ws.Range('C33') 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
Set-TextValue $ws.Range('D33') '0.999'
Set-TextValue $ws.Range('E33') '  +0.05%  '
Set-TextValue $ws.Range('B34') 'Filecoin'
Set-TextValue $ws.Range('C34') 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
Set-TextValue $ws.Range('D34') '5.86'
Set-TextValue $ws.Range('E34') '  -0.73%  '
Set-TextValue $ws.Range('B35') 'Mantle'
Set-TextValue $ws.Range('C35') 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
Set-TextValue $ws.Range('D35') '0.978'
Set-TextValue $ws.Range('E35') '  -2.65%  '
Set-TextValue $ws.Range('B36') 'Arweave'
Set-TextValue $ws.Range('C36') 'https://coinranking.com/coin/7XWg41D1+arweave-ar'
Set-TextValue $ws.Range('D36') '47.26'
Set-TextValue $ws.Range('E36') '  -1.14%  '
Set-TextValue $ws.Range('B37') 'OKB'
Set-TextValue $ws.Range('C37') 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
Set-TextValue $ws.Range('D37') '50.06'
Set-TextValue $ws.Range('E37') '  -0.32%  '
Set-TextValue $ws.Range('B38') 'Stacks'
Set-TextValue $ws.Range('C38') 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
Set-TextValue $ws.Range('D38') '2.05'
Set-TextValue $ws.Range('E38') '  -3.91%  '
Set-TextValue $ws.Range('B39') 'TheGraph'
Set-TextValue $ws.Range('C39') 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
Set-TextValue $ws.Range('D39') '0.312'
Set-TextValue $ws.Range('E39') '  -1.40%  '
Set-TextValue $ws.Range('B40') 'Kaspa'
Set-TextValue $ws.Range('C40') 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
Set-TextValue $ws.Range('D40') '0.124'
Set-TextValue $ws.Range('E40') '  +1.44%  '
Set-TextValue $ws.Range('B41') 'Cosmos'
Set-TextValue $ws.Range('C41') 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
Set-TextValue $ws.Range('D41') '8.55'
Set-TextValue $ws.Range('E41') '  -1.50%  '
Set-TextValue $ws.Range('B42') 'Maker'
Set-TextValue $ws.Range('C42') 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
Set-TextValue $ws.Range('D42') '2.816.13'
Set-TextValue $ws.Range('E42') '  -0.36%  '
Set-TextValue $ws.Range('D43') '384.32'
Set-TextValue $ws.Range('E43') '  +0.25%  '
Set-TextValue $ws.Range('B44') 'dogwifhat'
Set-TextValue $ws.Range('C44') 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
Set-TextValue $ws.Range('D44') '2.58'
Set-TextValue $ws.Range('E44') '  -8.12%  '
Set-TextValue $ws.Range('B45') 'VeChain'
Set-TextValue $ws.Range('C45') 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
Set-TextValue $ws.Range('D45') '0.0352'
Set-TextValue $ws.Range('E45') '  -2.14%  '
Set-TextValue $ws.Range('B46') 'Monero'
Set-TextValue $ws.Range('C46') 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
Set-TextValue $ws.Range('D46') '135.51'
Set-TextValue $ws.Range('E46') '  +0.47%  '
Set-TextValue $ws.Range('B47') 'USDe'
Set-TextValue $ws.Range('C47') 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
Set-TextValue $ws.Range('D47') '1.00'
Set-TextValue $ws.Range('E47') '  +0.01%  '
Set-TextValue $ws.Range('B48') 'InjectiveProtocol'
Set-TextValue $ws.Range('C48') 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
Set-TextValue $ws.Range('D48') '25.05'
Set-TextValue $ws.Range('E48') '  +0.51%  '
Set-TextValue $ws.Range('B49') 'ThetaToken'
Set-TextValue $ws.Range('C49') 'https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta'
Set-TextValue $ws.Range('D49') '2.22'
Set-TextValue $ws.Range('E49') '  -0.79%  '
Set-TextValue $ws.Range('B50') 'Stellar'
Set-TextValue $ws.Range('C50') 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
Set-TextValue $ws.Range('D50') '0.108'
Set-TextValue $ws.Range('E50') '  -0.60%  '
Set-TextValue $ws.Range('B51') 'THORChain'
Set-TextValue $ws.Range('C51') 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
Set-TextValue $ws.Range('D51') '6.78'
Set-TextValue $ws.Range('E51') '  -1.48%  '
